$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# --- Sheet "Produtos": remove the Shampoo/Gato row, finish the Ração row ---
$wsProdutos = $wb.Worksheets.Item("Produtos")
$wsProdutos.Rows(3).Delete()
$wsProdutos.Range("E2").Value = "Pacote"

$wsProdutos.Range("D2").NumberFormat = "@"
$wsProdutos.Range("D2").Value = "170"
$wsProdutos.Range("Z100").Copy()
$wsProdutos.Range("D2").PasteSpecial($xlPasteFormats)

# --- Sheet "Estoque": collapse the three rows into a single summary row ---
$wsEstoque = $wb.Worksheets.Item("Estoque")
$wsEstoque.Rows(4).Delete()
$wsEstoque.Rows(2).Delete()

$wsEstoque.Range("C2").NumberFormat = "@"
$wsEstoque.Range("C2").Value = "3"
$wsEstoque.Range("Z100").Copy()
$wsEstoque.Range("C2").PasteSpecial($xlPasteFormats)

$wsEstoque.Range("D2").NumberFormat = "@"
$wsEstoque.Range("D2").Value = "510.0"
$wsEstoque.Range("Z100").Copy()
$wsEstoque.Range("D2").PasteSpecial($xlPasteFormats)

# --- Sheet "Vendas": clear out the old sample row (start of the Sell window) ---
$wsVendas = $wb.Worksheets.Item("Vendas")
$wsVendas.Rows(2).Delete()

# --- New sheet "Métodos" at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMetodos = $wb.Worksheets.Add($null, $lastSheet)
$wsMetodos.Name = "Métodos"

$wsProdutos.Range("A1").Copy()
$wsMetodos.Range("B1").PasteSpecial($xlPasteFormats)
$wsMetodos.Range("B1").Value = "Métodos"

$wsProdutos.Range("A1").Copy()
$wsMetodos.Range("A2").PasteSpecial($xlPasteFormats)
$wsMetodos.Range("A2").Value = 0

$wsMetodos.Range("B2").Value = "Pacote"

$wsProdutos.Activate()
